# Update price/volume/ticker data so every row reflects the CRWD (CrowdStrike)
# series instead of the various other tickers that had been mistakenly mixed
# into this sheet ("fixed issues resulted from extra files").
#
# For every data row (2..32) we overwrite:
#   D = open_price, E = close_price, F = high_price, G = low_price,
#   H = shares_outstanding, I = fixed_ticker (-> "CRWD")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 63.5, 68.29000091552734, 79.79000091552734, 56, 250955140),
    @(3, 63.5, 68.29000091552734, 79.79000091552734, 56, 250955140),
    @(4, 63.5, 68.29000091552734, 79.79000091552734, 56, 250955140),
    @(5, 63.5, 68.29000091552734, 79.79000091552734, 56, 250955140),
    @(6, 63.5, 68.29000091552734, 79.79000091552734, 56, 250955140),
    @(7, 63.5, 68.29000091552734, 79.79000091552734, 56, 250955140),
    @(8, 90.06999969482422, 81.27999877929688, 101.879997253418, 78.88999938964844, 250955140),
    @(9, 49.93999862670898, 58, 60, 44.79999923706055, 250955140),
    @(10, 61.13000106811523, 59.63999938964844, 67.55999755859375, 54.34500122070312, 250955140),
    @(11, 66.27999877929688, 87.80999755859375, 88, 65.80000305175781, 250955140),
    @(12, 114.4800033569336, 125.7300033569336, 126.5500030517578, 93.37000274658205, 250955140),
    @(13, 124.5500030517578, 153.2799987792969, 155, 118.0999984741211, 250955140),
    @(14, 216.5599975585937, 216, 251.2799987792969, 198.1620025634766, 250955140),
    @(15, 209, 222.1499938964844, 225.9250030517578, 177.2870025634766, 250955140),
    @(16, 254.259994506836, 281, 289.239990234375, 227.25, 250955140),
    @(17, 275.5409851074219, 217.1399993896484, 298.4800109863281, 216.4900054931641, 250955140),
    @(18, 182.009994506836, 195.2100067138672, 197.25, 153.5, 250955140),
    @(19, 198.5, 159.9900054931641, 207.8300018310547, 130, 250955140),
    @(20, 180.7799987792969, 182.6100006103516, 205.729995727539, 178.3849945068359, 250955140),
    @(21, 165, 117.6500015258789, 168.4799957275391, 108.8899993896484, 250955140),
    @(22, 105.9000015258789, 120.6900024414062, 123.2399978637695, 105.0999984741211, 250955140),
    @(23, 118.5599975585938, 160.1300048828125, 161.7949981689453, 115.6699981689453, 250955140),
    @(24, 161, 163.0299987792969, 166.9900054931641, 140.5200042724609, 250955140),
    @(25, 178.1199951171875, 236.9900054931641, 238.479995727539, 174.8099975585938, 250955140),
    @(26, 296.1400146484375, 324.1499938964844, 338.4530029296875, 273.1600036621094, 250955140),
    @(27, 290.4400024414062, 313.6700134277344, 358.8399963378906, 288.5499877929688, 250955140),
    @(28, 232, 277.2799987792969, 285.7300109863281, 200.8099975585937, 250955140),
    @(29, 296.7900085449219, 345.9700012207031, 377.4200134277344, 295.0799865722656, 250955140),
    @(30, 386.9450073242188, 389.6600036621094, 455.5899963378906, 369.3099975585938, 250955140),
    @(31, 431.510009765625, 471.3699951171875, 474.2300109863281, 404.6300048828125, 250955140),
    @(32, 448.9500122070313, 423.7000122070313, 457.7999877929688, 409.25, 250955140)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 5).Value = $row[2]
    $ws.Cells.Item($r, 6).Value = $row[3]
    $ws.Cells.Item($r, 7).Value = $row[4]
    $ws.Cells.Item($r, 8).Value = $row[5]
    $ws.Cells.Item($r, 9).Value = "CRWD"
}
